$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 22: section title (merged A22:G22) + new parameter headers
$ws.Range("A22").Value = "Вывод: Холодильник не нужен, так как t'2 < t2 и по расчетам, я получил  N = 1071Вт, значит округляем до 1100 Вт"
$ws.Range("A22:G22").HorizontalAlignment = -4108
$ws.Range("A22:G22").Merge()

$ws.Range("H22").Value = "N(Вт)"
$ws.Range("I22").Value = "Nтабл(Вт)"
$ws.Range("J22").Value = "l3(м)"
$ws.Range("K22").Value = "m2(кг)"
$ws.Range("L22").Value = "t2(°C)"
$ws.Range("M22").Value = "t'2(°C)"

# Row 23: formulas/values computed from row 22 headers
$ws.Range("H23").Formula = "=N5*O4*((D12/Q5)+(D12/(2*P5)))"
$ws.Range("I23").Value = 1100
$ws.Range("J23").Formula = "=S4/(D11*((1/D16)-1))"
$ws.Range("K23").Formula = "=1000*((M4*(POWER(K5,2)-POWER(J5,2)))/4)*D17"
$ws.Range("L23").Formula = "=33.2+((O4/D18)*(2*0.2-0.05))"
$ws.Range("M23").Formula = "=(0.2*D12)+23+(E13/(D16*S4))"

# Match formatting used by the rest of the table (style index 2: centered)
$ws.Range("H22:M23").VerticalAlignment = -4108
$ws.Range("H22:M23").HorizontalAlignment = -4108

# Sheet view tweaks recorded in the diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A23").Select()

Write-Host "done"
